$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fix the garbled worker name in D17:D21 so it matches the clean
# "MALENA PATRICIA DIAZ CAÑATE" already used in D16.
$ws.Range("D17").Value = "MALENA PATRICIA DIAZ CAÑATE"
$ws.Range("D18").Value = "MALENA PATRICIA DIAZ CAÑATE"
$ws.Range("D19").Value = "MALENA PATRICIA DIAZ CAÑATE"
$ws.Range("D20").Value = "MALENA PATRICIA DIAZ CAÑATE"
$ws.Range("D21").Value = "MALENA PATRICIA DIAZ CAÑATE"

# Reorder the "Periodo Mora" column from descending to ascending:
# before = 2411,2410,2409,2408,2407,2406 (rows 16..21)
# after  = 2406,2407,2408,2409,2410,2411 (rows 16..21)
$ws.Range("E16").Value = "2406"
$ws.Range("E17").Value = "2407"
$ws.Range("E18").Value = "2408"
$ws.Range("E19").Value = "2409"
$ws.Range("E20").Value = "2410"
$ws.Range("E21").Value = "2411"
